$d = $word.ActiveDocument

$replacements = @(
    @{old="266÷7="; new="261÷5="},
    @{old="835÷4="; new="225÷3="},
    @{old="191÷9="; new="110÷6="},
    @{old="854÷4="; new="186÷6="},
    @{old="426÷8="; new="902÷6="},
    @{old="731÷4="; new="299÷3="},
    @{old="157÷5="; new="788÷6="},
    @{old="998÷7="; new="895÷2="},
    @{old="893÷6="; new="989÷5="},
    @{old="365÷9="; new="933÷8="},
    @{old="469÷3="; new="586÷9="},
    @{old="560÷4="; new="424÷2="},
    @{old="621÷4="; new="850÷6="},
    @{old="112÷8="; new="917÷8="},
    @{old="719÷5="; new="490÷3="},
    @{old="224÷9="; new="475÷7="},
    @{old="466÷6="; new="642÷9="},
    @{old="838÷7="; new="787÷6="},
    @{old="154÷6="; new="427÷6="},
    @{old="421÷7="; new="260÷7="},
    @{old="659÷9="; new="643÷8="},
    @{old="124÷5="; new="491÷5="},
    @{old="153÷2="; new="831÷4="},
    @{old="119÷4="; new="722÷8="},
    @{old="142÷4="; new="878÷9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
